$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1251.62255859375
$ws.Range("C2").Value = 0.9389999999999999
$ws.Range("D2").Value = 0.9298999905586243
$ws.Range("E2").Value = 1.380499958992004
$ws.Range("F2").Value = 0.6366000175476074
$ws.Range("H2").Value = 0.9399

$ws.Range("B3").Value = 1225.838745117188
$ws.Range("C3").Value = 0.9807
$ws.Range("D3").Value = 0.948
$ws.Range("E3").Value = 1.483100056648254
$ws.Range("F3").Value = 0.6736000180244446
$ws.Range("H3").Value = 1.1008

$ws.Range("B4").Value = 811.7313842773438
$ws.Range("C4").Value = 0.9483
$ws.Range("D4").Value = 0.905
$ws.Range("E4").Value = 1.615599989891052
$ws.Range("F4").Value = 0.7594000101089478
$ws.Range("H4").Value = 0.7193000000000001

$ws.Range("B5").Value = 799.4152221679688
$ws.Range("C5").Value = 0.8433
$ws.Range("D5").Value = 0.8367
$ws.Range("E5").Value = 1.105599999427795
$ws.Range("F5").Value = 0.5630000233650208
$ws.Range("H5").Value = 0.1147

$ws.Range("B6").Value = 1101.406494140625
$ws.Range("C6").Value = 0.8714
$ws.Range("D6").Value = 0.8665
$ws.Range("E6").Value = 1.078999996185303
$ws.Range("F6").Value = 0.679099977016449
$ws.Range("H6").Value = 0.3784

$ws.Range("B7").Value = 855.4713745117188
$ws.Range("C7").Value = 0.8615
$ws.Range("D7").Value = 0.862500011920929
$ws.Range("E7").Value = 1.006800055503845
$ws.Range("F7").Value = 0.7088000178337097
$ws.Range("H7").Value = 0.3429

$ws.Range("B8").Value = 940.2319946289062
$ws.Range("C8").Value = 0.8425
$ws.Range("D8").Value = 0.8421
$ws.Range("E8").Value = 1.075199961662292
$ws.Range("F8").Value = 0.7229999899864197
$ws.Range("H8").Value = 0.1623

$ws.Range("B9").Value = 6985.7177734375
$ws.Range("C9").Value = 0.9003
$ws.Range("D9").Value = 0.8848
$ws.Range("E9").Value = 1.615599989891052
$ws.Range("F9").Value = 0.5630000233650208
$ws.Range("H9").Value = 3.7583
